# Updated cryptos list values (price & 1h volume change) to match latest
# scrape results. Price column (D) must stay text (some values like
# "28.606.35" or "1.00" are not meant to be interpreted as numbers), so we
# force the text number-format before assigning.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.606.35"
$ws.Range("E2").Value = "  +1.06%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.563.34"
$ws.Range("E3").Value = "  +0.21%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "210.11"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.508"
$ws.Range("E6").Value = "  +3.80%  "

$ws.Range("E7").Value = "  -0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "24.87"
$ws.Range("E8").Value = "  +5.61%  "

$ws.Range("E9").Value = "  +0.84%  "

$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.787.14"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.566.69"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.654.21"
$ws.Range("E14").Value = "  +1.24%  "

$ws.Range("E15").Value = "  +0.98%  "

$ws.Range("E16").Value = "  -0.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.41"
$ws.Range("E17").Value = "  +0.74%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.87"
$ws.Range("E18").Value = "  -0.19%  "

$ws.Range("E19").Value = "  -0.51%  "

$ws.Range("E20").Value = "  +0.38%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  -0.31%  "

$ws.Range("E22").Value = "  -0.11%  "

$ws.Range("E24").Value = "  +1.10%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.47"
$ws.Range("E25").Value = "  +0.92%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.106"
$ws.Range("E26").Value = "  +2.73%  "

$ws.Range("E27").Value = "  -0.71%  "

$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.22"
$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("E30").Value = "  -3.81%  "

$ws.Range("E31").Value = "  -0.39%  "

$ws.Range("E32").Value = "  +0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.399.01"
$ws.Range("E33").Value = "  +1.70%  "

$ws.Range("E34").Value = "  -2.91%  "

$ws.Range("E35").Value = "  -2.51%  "

$ws.Range("E36").Value = "  -1.76%  "

$ws.Range("E37").Value = "  +1.32%  "

$ws.Range("E38").Value = "  -2.04%  "

$ws.Range("E39").Value = "  -0.34%  "

$ws.Range("E41").Value = "  -0.42%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  -0.13%  "

$ws.Range("E43").Value = "  -1.76%  "

$ws.Range("E44").Value = "  -2.46%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "63.84"
$ws.Range("E45").Value = "  +2.92%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.699.64"
$ws.Range("E47").Value = "  +0.27%  "

$ws.Range("E48").Value = "  -5.33%  "

$ws.Range("E49").Value = "  -0.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "42.21"
$ws.Range("E50").Value = "  +4.94%  "

$ws.Range("E51").Value = "  -0.66%  "

